# culture_collection を MIxS から再度削除
# Remove the "culture_collection" column (column V) entirely: its header cell
# value, together with every column to its right, shifts one position to the
# left. Column comments are a separate object model from cell values, so they
# have to be re-aligned by hand after the column delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRow = 15
$deleteCol = 22       # V = culture_collection
$firstShiftCol = 23   # W, first column to the right of the deleted one
$lastCol = 72         # BT, last used column in the header row (before delete)

# 1) Delete the culture_collection column. This shifts every cell value one
#    column to the left (and keeps row "spans" / sheet dimension correct),
#    but comments are NOT moved by Excel when a column is deleted this way.
$ws.Columns("V:V").Delete()

# 2) Re-align the comments to match: comment that was on column c (W..BT)
#    belongs on column c-1 (V..BS) now.
for ($c = $firstShiftCol; $c -le $lastCol; $c++) {
    $srcCell = $ws.Cells.Item($headerRow, $c)
    $dstCell = $ws.Cells.Item($headerRow, $c - 1)

    $srcComment = $srcCell.Comment
    if ($srcComment) {
        $commentText = $srcComment.Text()

        $dstComment = $dstCell.Comment
        if ($dstComment) {
            $dstComment.Text($commentText)
        } else {
            $dstCell.AddComment($commentText)
        }
    }
}

# 3) The comment left behind on the old last column (BT) is now a duplicate
#    of BS's comment; remove it.
$lastComment = $ws.Cells.Item($headerRow, $lastCol).Comment
if ($lastComment) {
    $lastComment.Delete()
}
